$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Card10")

$row = 14

# Columns A and L hold text that Excel would otherwise auto-convert
# (A14 looks like a plain number, L14 looks like a date) -- force them
# to be stored as text by pre-setting the cell number format to "@".
$ws.Cells.Item($row, 1).NumberFormat = "@"
$ws.Cells.Item($row, 1).Value = "10"

$ws.Cells.Item($row, 12).NumberFormat = "@"
$ws.Cells.Item($row, 12).Value = "5/3/2025"

$ws.Cells.Item($row, 13).Value = "قطع سير كويلر مسنن 1270"
$ws.Cells.Item($row, 14).Value = "تم تغير سير 1270"
$ws.Cells.Item($row, 15).Value = "فني"
